$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "***"
$ws.Range("D2").Value = "**"
$ws.Range("E2").Value = "*"
$ws.Range("H2").Value = " "
$ws.Range("L2").Value = "*"
$ws.Range("M2").Value = "*"
$ws.Range("H3").Value = "*"
$ws.Range("M3").Value = " "
$ws.Range("M4").Value = "*"
$ws.Range("B5").Value = "***"
$ws.Range("D5").Value = "*"
$ws.Range("H5").Value = "*"
$ws.Range("L5").Value = "**"
$ws.Range("M5").Value = "***"
$ws.Range("D6").Value = "*"
$ws.Range("M6").Value = "*"
$ws.Range("B7").Value = "**"
$ws.Range("E7").Value = "*"
$ws.Range("H8").Value = " "
$ws.Range("M8").Value = "*"
$ws.Range("D9").Value = "*"
$ws.Range("E9").Value = "***"
$ws.Range("G9").Value = "***"
$ws.Range("D10").Value = "***"
$ws.Range("E10").Value = "***"
$ws.Range("G10").Value = "**"
$ws.Range("M10").Value = "*"
$ws.Range("B11").Value = "*"
$ws.Range("H11").Value = " "
$ws.Range("L11").Value = "***"
$ws.Range("M11").Value = " "
$ws.Range("B12").Value = " "
$ws.Range("D12").Value = "*"
$ws.Range("E12").Value = " "
$ws.Range("H12").Value = "**"
$ws.Range("L12").Value = "*"
$ws.Range("D13").Value = "**"
$ws.Range("E13").Value = " "
$ws.Range("M13").Value = "*"
$ws.Range("D14").Value = "***"
$ws.Range("E14").Value = " "
$ws.Range("E15").Value = "*"
$ws.Range("K15").Value = "***"
$ws.Range("B16").Value = "*"
$ws.Range("E16").Value = "*"
$ws.Range("H16").Value = " "
$ws.Range("M16").Value = "*"
$ws.Range("C17").Value = " "
$ws.Range("D17").Value = "***"
$ws.Range("D18").Value = "**"
$ws.Range("D19").Value = "***"
$ws.Range("E19").Value = "*"
$ws.Range("H19").Value = " "
$ws.Range("M19").Value = "***"
$ws.Range("H20").Value = "*"
$ws.Range("C21").Value = "*"
$ws.Range("F21").Value = "**"
$ws.Range("C22").Value = "***"
$ws.Range("G22").Value = "***"
$ws.Range("I22").Value = "**"
$ws.Range("J22").Value = "***"
$ws.Range("B23").Value = "*"
$ws.Range("H23").Value = "***"
$ws.Range("I23").Value = "**"
